# Auto-generated script applying scheduled-runner price/profit updates
# to the Kujata_Profits workbook across multiple job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4485
$ws.Range("I64").Value = 4983.3335
$ws.Range("K64").Value = 4983.3335
$ws.Range("M64").Value = -4735.3335
$ws.Range("H67").Value = 4485
$ws.Range("I67").Value = 4983.3335
$ws.Range("K67").Value = 4983.3335
$ws.Range("M67").Value = -4125.3335
$ws.Range("H116").Value = 2924.9412
$ws.Range("I116").Value = 2229.2856
$ws.Range("J116").Value = 3411.9
$ws.Range("K116").Value = 2229.2856
$ws.Range("L116").Value = 3411.9
$ws.Range("M116").Value = 1212.7144
$ws.Range("N116").Value = -10295.9
$ws.Range("H129").Value = 848.8333
$ws.Range("I129").Value = 301.36365
$ws.Range("J129").Value = 988.8837
$ws.Range("K129").Value = 904.09095
$ws.Range("L129").Value = 2966.6511
$ws.Range("M129").Value = 4095.90905
$ws.Range("N129").Value = -12966.6511
$ws.Range("H132").Value = 4834526
$ws.Range("I132").Value = 6063388.5
$ws.Range("J132").Value = 6851.5
$ws.Range("K132").Value = 18190165.5
$ws.Range("L132").Value = 20554.5
$ws.Range("M132").Value = -18187635.5
$ws.Range("N132").Value = -25614.5
$ws.Range("H135").Value = 752.4167
$ws.Range("J135").Value = 1968
$ws.Range("L135").Value = 17712
$ws.Range("N135").Value = -22782
$ws.Range("H137").Value = 1100.2333
$ws.Range("I137").Value = 822
$ws.Range("J137").Value = 1378.4667
$ws.Range("K137").Value = 2466
$ws.Range("L137").Value = 4135.4001
$ws.Range("M137").Value = 84
$ws.Range("N137").Value = -9235.400099999999
$ws.Range("H138").Value = 1485.07
$ws.Range("I138").Value = 793.6
$ws.Range("J138").Value = 1607.0941
$ws.Range("K138").Value = 2380.8
$ws.Range("L138").Value = 4821.2823
$ws.Range("M138").Value = 2759.2
$ws.Range("N138").Value = -15101.2823

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2172.2903
$ws.Range("I132").Value = 2025.6471
$ws.Range("J132").Value = 2350.3572
$ws.Range("K132").Value = 6076.9413
$ws.Range("L132").Value = 7051.071599999999
$ws.Range("M132").Value = -3546.9413
$ws.Range("N132").Value = -12111.0716

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2849.0645
$ws.Range("I86").Value = 3196.261
$ws.Range("J86").Value = 1850.875
$ws.Range("K86").Value = 3196.261
$ws.Range("L86").Value = 1850.875
$ws.Range("M86").Value = -2073.261
$ws.Range("N86").Value = -4096.875
$ws.Range("H87").Value = 53000
$ws.Range("J87").Value = 53000
$ws.Range("L87").Value = 53000
$ws.Range("N87").Value = -55496
$ws.Range("H89").Value = 2849.0645
$ws.Range("I89").Value = 3196.261
$ws.Range("J89").Value = 1850.875
$ws.Range("K89").Value = 15981.305
$ws.Range("L89").Value = 9254.375
$ws.Range("M89").Value = -10365.305
$ws.Range("N89").Value = -20486.375
$ws.Range("H90").Value = 53000
$ws.Range("J90").Value = 53000
$ws.Range("L90").Value = 159000
$ws.Range("N90").Value = -171480
$ws.Range("H134").Value = 3688.6738
$ws.Range("I134").Value = 932.55554
$ws.Range("K134").Value = 2797.66662
$ws.Range("M134").Value = -262.66662

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 23257.334
$ws.Range("J50").Value = 23257.334
$ws.Range("L50").Value = 23257.334
$ws.Range("N50").Value = -24507.334
$ws.Range("H60").Value = 6124.75
$ws.Range("I60").Value = 3142.8572
$ws.Range("J60").Value = 26998
$ws.Range("K60").Value = 3142.8572
$ws.Range("L60").Value = 26998
$ws.Range("M60").Value = -2631.8572
$ws.Range("N60").Value = -28020
$ws.Range("H62").Value = 200000000
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 200000000
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H68").Value = 16500
$ws.Range("J68").Value = 16500
$ws.Range("L68").Value = 16500
$ws.Range("N68").Value = -17998
$ws.Range("H71").Value = 16500
$ws.Range("J71").Value = 16500
$ws.Range("L71").Value = 49500
$ws.Range("N71").Value = -56988
$ws.Range("H132").Value = 3246.842
$ws.Range("I132").Value = 3487.422
$ws.Range("J132").Value = 2344.6667
$ws.Range("K132").Value = 10462.266
$ws.Range("L132").Value = 7034.000100000001
$ws.Range("M132").Value = -7932.266
$ws.Range("N132").Value = -12094.0001
$ws.Range("H134").Value = 1007.1042
$ws.Range("I134").Value = 1065.6842
$ws.Range("J134").Value = 784.5
$ws.Range("K134").Value = 3197.0526
$ws.Range("L134").Value = 2353.5
$ws.Range("M134").Value = -662.0526
$ws.Range("N134").Value = -7423.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 70.2
$ws.Range("I2").Value = 37.75
$ws.Range("K2").Value = 226.5
$ws.Range("M2").Value = -113.5
$ws.Range("H34").Value = 1860.25
$ws.Range("I34").Value = 1297.3334
$ws.Range("J34").Value = 2198
$ws.Range("K34").Value = 3892.0002
$ws.Range("L34").Value = 6594
$ws.Range("M34").Value = -3808.0002
$ws.Range("N34").Value = -6762
$ws.Range("H39").Value = 3142.6155
$ws.Range("J39").Value = 2945.4
$ws.Range("L39").Value = 8836.200000000001
$ws.Range("N39").Value = -9424.200000000001
$ws.Range("H55").Value = 2668
$ws.Range("I55").Value = 1004
$ws.Range("J55").Value = 3500
$ws.Range("K55").Value = 3012
$ws.Range("L55").Value = 10500
$ws.Range("M55").Value = -2835
$ws.Range("N55").Value = -10854
$ws.Range("H104").Value = 3905
$ws.Range("I104").Value = 3325.6667
$ws.Range("J104").Value = 4194.6665
$ws.Range("K104").Value = 9977.000100000001
$ws.Range("L104").Value = 12583.9995
$ws.Range("M104").Value = -7356.000100000001
$ws.Range("N104").Value = -17825.9995
$ws.Range("H131").Value = 22728572
$ws.Range("I131").Value = 100000264
$ws.Range("J131").Value = 1603.5883
$ws.Range("K131").Value = 300000792
$ws.Range("L131").Value = 4810.7649
$ws.Range("M131").Value = -299995752
$ws.Range("N131").Value = -14890.7649
$ws.Range("H140").Value = 23102.408
$ws.Range("I140").Value = 52425.41
$ws.Range("J140").Value = 2942.8438
$ws.Range("K140").Value = 157276.23
$ws.Range("L140").Value = 8828.5314
$ws.Range("M140").Value = -152096.23
$ws.Range("N140").Value = -19188.5314

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 32145016
$ws.Range("I70").Value = 35716344
$ws.Range("J70").Value = 28573686
$ws.Range("K70").Value = 35716344
$ws.Range("L70").Value = 28573686
$ws.Range("M70").Value = -35716074
$ws.Range("N70").Value = -28574226
$ws.Range("H73").Value = 32145016
$ws.Range("I73").Value = 35716344
$ws.Range("J73").Value = 28573686
$ws.Range("K73").Value = 35716344
$ws.Range("L73").Value = 28573686
$ws.Range("M73").Value = -35715408
$ws.Range("N73").Value = -28575558
$ws.Range("H132").Value = 2604.8667
$ws.Range("I132").Value = 2255.0833
$ws.Range("J132").Value = 4004
$ws.Range("K132").Value = 6765.249899999999
$ws.Range("L132").Value = 12012
$ws.Range("M132").Value = -4235.249899999999
$ws.Range("N132").Value = -17072

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1210.2667
$ws.Range("I7").Value = 1034.5385
$ws.Range("J7").Value = 2352.5
$ws.Range("K7").Value = 1034.5385
$ws.Range("L7").Value = 2352.5
$ws.Range("M7").Value = -922.5385000000001
$ws.Range("N7").Value = -2576.5
$ws.Range("H46").Value = 4791.6665
$ws.Range("I46").Value = 1333.3334
$ws.Range("J46").Value = 5944.4443
$ws.Range("K46").Value = 1333.3334
$ws.Range("L46").Value = 5944.4443
$ws.Range("M46").Value = -1145.3334
$ws.Range("N46").Value = -6320.4443
$ws.Range("H61").Value = 1328.75
$ws.Range("I61").Value = 1103.3334
$ws.Range("K61").Value = 1103.3334
$ws.Range("M61").Value = -901.3334
$ws.Range("H70").Value = 28666.334
$ws.Range("J70").Value = 28666.334
$ws.Range("L70").Value = 28666.334
$ws.Range("N70").Value = -29206.334
$ws.Range("H73").Value = 28666.334
$ws.Range("J73").Value = 28666.334
$ws.Range("L73").Value = 28666.334
$ws.Range("N73").Value = -30538.334
$ws.Range("H113").Value = 1328.75
$ws.Range("I113").Value = 1103.3334
$ws.Range("K113").Value = 1103.3334
$ws.Range("M113").Value = 1066.6666
$ws.Range("H126").Value = 1210.2667
$ws.Range("I126").Value = 1034.5385
$ws.Range("J126").Value = 2352.5
$ws.Range("K126").Value = 3103.6155
$ws.Range("L126").Value = 7057.5
$ws.Range("M126").Value = -633.6155000000003
$ws.Range("N126").Value = -11997.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 5000
$ws.Range("J10").Value = 5000
$ws.Range("L10").Value = 5000
$ws.Range("N10").Value = -5338
$ws.Range("H113").Value = 419.3125
$ws.Range("I113").Value = 310.5
$ws.Range("J113").Value = 528.125
$ws.Range("K113").Value = 931.5
$ws.Range("L113").Value = 1584.375
$ws.Range("M113").Value = 1238.5
$ws.Range("N113").Value = -5924.375
